# Monster Hunter General Data - "did time betweem releases"
#
# The release table was missing a row for "Monster Hunter X" (JPN, 3DS,
# 2015-11-28, Generation 4, Director Yasunori Ichinose). That row is
# inserted at row 21 (pushing every following row down by one), and the
# surrounding rows - which had drifted out of sync with the correct
# title/region/date/director/console combinations - are corrected back
# into alignment. A trailing row for "Monster Hunter Rise" (JPN) is
# (re)appended at the end so both the NA and JPN Rise releases are listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21; everything from the old row 21 onward shifts
# down to make room (old row 21 "Generations" becomes row 22, etc.)
$ws.Rows.Item(21).Insert()

# New row 21: Monster Hunter X
$ws.Cells.Item(21, 1).Value = "Monster Hunter X"
$ws.Cells.Item(21, 2).Value = "JPN"
$ws.Cells.Item(21, 3).Value = 42336
$ws.Cells.Item(21, 4).Value = 4
$ws.Cells.Item(21, 5).Value = "Yasunori Ichinose"
$ws.Cells.Item(21, 6).Value = "3DS"

# Row 22: Monster Hunter XX (unchanged data, just shifted down)
$ws.Cells.Item(22, 1).Value = "Monster Hunter XX"
$ws.Cells.Item(22, 2).Value = "JPN"
$ws.Cells.Item(22, 3).Value = 42812
$ws.Cells.Item(22, 4).Value = 4
$ws.Cells.Item(22, 5).Value = "Yasunori Ichinose"
$ws.Cells.Item(22, 6).Value = "3DS"

# Row 23: Monster Hunter Generations, NA
$ws.Cells.Item(23, 1).Value = "Monster Hunter Generations"
$ws.Cells.Item(23, 2).Value = "NA"
$ws.Cells.Item(23, 3).Value = 42566
$ws.Cells.Item(23, 4).Value = 4
$ws.Cells.Item(23, 5).Value = "Yasunori Ichinose"
$ws.Cells.Item(23, 6).Value = "3DS"

# Row 24: Monster Hunter Generations Ultimate, JPN
$ws.Cells.Item(24, 1).Value = "Monster Hunter Generations Ultimate"
$ws.Cells.Item(24, 2).Value = "JPN"
$ws.Cells.Item(24, 3).Value = 42972
$ws.Cells.Item(24, 4).Value = 4
$ws.Cells.Item(24, 5).Value = "Yasunori Ichinose"
$ws.Cells.Item(24, 6).Value = "Nintendo Switch"

# Row 25: Monster Hunter Generations Ultimate, NA
$ws.Cells.Item(25, 1).Value = "Monster Hunter Generations Ultimate"
$ws.Cells.Item(25, 2).Value = "NA"
$ws.Cells.Item(25, 3).Value = 43340
$ws.Cells.Item(25, 4).Value = 4
$ws.Cells.Item(25, 5).Value = "Yasunori Ichinose"
$ws.Cells.Item(25, 6).Value = "Nintendo Switch"

# Row 26: Monster Hunter: World, JPN
$ws.Cells.Item(26, 1).Value = "Monster Hunter: World"
$ws.Cells.Item(26, 2).Value = "JPN"
$ws.Cells.Item(26, 3).Value = 43126
$ws.Cells.Item(26, 4).Value = 5
$ws.Cells.Item(26, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(26, 6).Value = "Xbox One"

# Row 27: Monster Hunter: World, NA
$ws.Cells.Item(27, 1).Value = "Monster Hunter: World"
$ws.Cells.Item(27, 2).Value = "NA"
$ws.Cells.Item(27, 3).Value = 43126
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(27, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(27, 6).Value = "PlayStation 4"

# Row 28: Monster Hunter: World, JPN
$ws.Cells.Item(28, 1).Value = "Monster Hunter: World"
$ws.Cells.Item(28, 2).Value = "JPN"
$ws.Cells.Item(28, 3).Value = 43126
$ws.Cells.Item(28, 4).Value = 5
$ws.Cells.Item(28, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(28, 6).Value = "Xbox One"

# Row 29: Monster Hunter: World, NA
$ws.Cells.Item(29, 1).Value = "Monster Hunter: World"
$ws.Cells.Item(29, 2).Value = "NA"
$ws.Cells.Item(29, 3).Value = 43126
$ws.Cells.Item(29, 4).Value = 5
$ws.Cells.Item(29, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(29, 6).Value = "PlayStation 4"

# Row 30: Monster Hunter World: Iceborne, JPN
$ws.Cells.Item(30, 1).Value = "Monster Hunter World: Iceborne"
$ws.Cells.Item(30, 2).Value = "JPN"
$ws.Cells.Item(30, 3).Value = 43714
$ws.Cells.Item(30, 4).Value = 5
$ws.Cells.Item(30, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(30, 6).Value = "Xbox One"

# Row 31: Monster Hunter World: Iceborne, NA
$ws.Cells.Item(31, 1).Value = "Monster Hunter World: Iceborne"
$ws.Cells.Item(31, 2).Value = "NA"
$ws.Cells.Item(31, 3).Value = 43714
$ws.Cells.Item(31, 4).Value = 5
$ws.Cells.Item(31, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(31, 6).Value = "PlayStation 4"

# Row 32: Monster Hunter World: Iceborne, JPN
$ws.Cells.Item(32, 1).Value = "Monster Hunter World: Iceborne"
$ws.Cells.Item(32, 2).Value = "JPN"
$ws.Cells.Item(32, 3).Value = 43714
$ws.Cells.Item(32, 4).Value = 5
$ws.Cells.Item(32, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(32, 6).Value = "Xbox One"

# Row 33: Monster Hunter World: Iceborne, NA
$ws.Cells.Item(33, 1).Value = "Monster Hunter World: Iceborne"
$ws.Cells.Item(33, 2).Value = "NA"
$ws.Cells.Item(33, 3).Value = 43714
$ws.Cells.Item(33, 4).Value = 5
$ws.Cells.Item(33, 5).Value = "Yuya Tokuda"
$ws.Cells.Item(33, 6).Value = "PlayStation 4"

# Row 34: Monster Hunter Rise, NA
$ws.Cells.Item(34, 1).Value = "Monster Hunter Rise"
$ws.Cells.Item(34, 2).Value = "NA"
$ws.Cells.Item(34, 3).Value = 44281
$ws.Cells.Item(34, 4).Value = 5
$ws.Cells.Item(34, 5).Value = "Yasunori Ichinose"
$ws.Cells.Item(34, 6).Value = "Nintendo Switch"

# Row 35 (new): Monster Hunter Rise, JPN
$ws.Cells.Item(35, 1).Value = "Monster Hunter Rise"
$ws.Cells.Item(35, 2).Value = "JPN"
$ws.Cells.Item(35, 3).Value = 44281
$ws.Cells.Item(35, 4).Value = 5
$ws.Cells.Item(35, 5).Value = "Yasunori Ichinose"
$ws.Cells.Item(35, 6).Value = "Nintendo Switch"

# Match the cursor position recorded in the saved file.
$ws.Range("C39").Select()
